{"js": "// The document stores per-plate metadata as inline pseudo-XML tags, e.g.\n//   <id>p081v_a1</id>\n// split across THREE runs: \"<id>\" (Courier New / gold), \"p081v_a1\" (plain),\n// \"</id>\" (Courier New / gold). The edit collapses each triple into a\n// SINGLE run \"<id>p081v_1</id>\" (dropping the \"a\" before the trailing\n// number) using the tag run's formatting.\n\nconst body = context.document.body;\n\n// Find every \"<id>\" opening tag and its matching \"</id>\" closing tag.\nconst idOpens = body.search(\"<id>\", { matchCase: true });\nconst idCloses = body.search(\"</id>\", { matchCase: true });\nidOpens.load(\"items\");\nidCloses.load(\"items\");\nawait context.sync();\n\nif (idOpens.items.length !== idCloses.items.length) {\n  throw new Error(\"Mismatched <id>/</id> counts: \" + idOpens.items.length + \" vs \" + idCloses.items.length);\n}\n\nfor (let i = 0; i < idOpens.items.length; i++) {\n  // The range spanning from \"<id>\" through \"</id>\" (inclusive) covers the\n  // three original runs: \"<id>\", \"p081v_aN\", \"</id>\".\n  const idRange = idOpens.items[i].expandTo(idCloses.items[i]);\n  idRange.load(\"text\");\n  await context.sync();\n\n  const match = idRange.text.match(/^<id>p081v_a(\\d+)<\\/id>$/);\n  if (!match) {\n    throw new Error(\"Unexpected <id> content: \" + idRange.text);\n  }\n\n  // Replacing the whole range with a single insertText call merges the\n  // three runs into one, picking up the formatting of the range's first\n  // run (the \"<id>\" tag's Courier New / gold styling).\n  idRange.insertText(\"<id>p081v_\" + match[1] + \"</id>\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document stores per-plate metadata as inline pseudo-XML tags, e.g.\n#   <id>p081v_a1</id>\n# split across THREE runs: \"<id>\" (Courier New / gold), \"p081v_a1\" (plain),\n# \"</id>\" (Courier New / gold). The edit collapses each triple into a\n# SINGLE run \"<id>p081v_1</id>\" (dropping the \"a\" before the trailing\n# number), using the tag run's formatting.\n\n$d = $word.ActiveDocument\n\n# Discover every \"<id>p081v_aN</id>\" occurrence from the plain document text\n# so we don't have to hard-code how many there are.\n$fullText = $d.Content.Text\n$idMatches = [regex]::Matches($fullText, '<id>p081v_a(\\d+)</id>')\n\nforeach ($m in $idMatches) {\n    $oldTag = $m.Value\n    $num = $m.Groups[1].Value\n    $newTag = \"<id>p081v_$num</id>\"\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldTag\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWildcards = $false\n    $found = $rng.Find.Execute()\n\n    if ($found) {\n        # Assigning .Text replaces the whole matched range (which spans the\n        # original three runs) with one run, inheriting the range's\n        # (i.e. the \"<id>\" run's) formatting.\n        $rng.Text = $newTag\n    }\n}\n"}
